$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CARGA")

$data = New-Object 'object[,]' 21,16
$data[0,0] = 'Matias Godoy'
$data[0,1] = 'Extremos'
$data[0,2] = 68
$data[0,3] = 3147.3613300000002
$data[0,4] = 46.370420000000003
$data[0,5] = 21.879809999999999
$data[0,6] = 83.68
$data[0,7] = 0
$data[0,8] = 13.86
$data[0,9] = 0
$data[0,10] = 127.51
$data[0,11] = 44.72
$data[0,12] = 0
$data[0,13] = 328.11998999999997
$data[0,14] = 46049
$data[0,15] = 5
$data[1,0] = 'Nicolas Thaller'
$data[1,1] = 'Defensor Central'
$data[1,2] = 54
$data[1,3] = 3804.2575099999999
$data[1,4] = 70.22251
$data[1,5] = 27.03379
$data[1,6] = 130.69
$data[1,7] = 6.92
$data[1,8] = 26.96
$data[1,9] = 1.69
$data[1,10] = 294.94999000000001
$data[1,11] = 187.35999000000001
$data[1,12] = 26.79
$data[1,13] = 618.82997999999998
$data[1,14] = 46049
$data[1,15] = 5
$data[2,0] = 'Lucas Cano'
$data[2,1] = 'Delantero Central'
$data[2,2] = 63
$data[2,3] = 4821.6980000000003
$data[2,4] = 76.142889999999994
$data[2,5] = 25.88083
$data[2,6] = 134.13999999999999
$data[2,7] = 3.23
$data[2,8] = 22.42
$data[2,9] = 0
$data[2,10] = 947.41002000000003
$data[2,11] = 189.14
$data[2,12] = 8.1199999999999992
$data[2,13] = 1033.8800100000001
$data[2,14] = 46049
$data[2,15] = 5
$data[3,0] = 'Pier Barrios'
$data[3,1] = 'Defensor Central'
$data[3,2] = 54
$data[3,3] = 3659.2968799999999
$data[3,4] = 67.546689999999998
$data[3,5] = 27.01426
$data[3,6] = 135.26
$data[3,7] = 0
$data[3,8] = 30.28
$data[3,9] = 0
$data[3,10] = 268.36
$data[3,11] = 61.74
$data[3,12] = 13.87
$data[3,13] = 574.45001000000002
$data[3,14] = 46049
$data[3,15] = 5
$data[4,0] = 'Federico Lertora'
$data[4,1] = 'Volante Central'
$data[4,2] = 54
$data[4,3] = 4099.0499
$data[4,4] = 75.664060000000006
$data[4,5] = 25.887599999999999
$data[4,6] = 137.52000000000001
$data[4,7] = 4.93
$data[4,8] = 39.9
$data[4,9] = 0.8
$data[4,10] = 328.94997999999998
$data[4,11] = 208.28
$data[4,12] = 13.07
$data[4,13] = 730.24000999999998
$data[4,14] = 46049
$data[4,15] = 5
$data[5,0] = 'Facundo Castro'
$data[5,1] = 'Delantero Central'
$data[5,2] = 68
$data[5,3] = 4137.7515000000003
$data[5,4] = 60.961950000000002
$data[5,5] = 26.045500000000001
$data[5,6] = 141.19
$data[5,7] = 1.19
$data[5,8] = 42.69
$data[5,9] = 0.9
$data[5,10] = 260.27999999999997
$data[5,11] = 136.75
$data[5,12] = 10.15
$data[5,13] = 708.56
$data[5,14] = 46049
$data[5,15] = 5
$data[6,0] = 'Zahir Ibarra'
$data[6,1] = 'Defensor Central'
$data[6,2] = 54
$data[6,3] = 4021.7503400000001
$data[6,4] = 74.237189999999998
$data[6,5] = 31.458629999999999
$data[6,6] = 166.75
$data[6,7] = 3.72
$data[6,8] = 58.85
$data[6,9] = 3.49
$data[6,10] = 314.77999999999997
$data[6,11] = 237.51
$data[6,12] = 165.84
$data[6,13] = 711.10001
$data[6,14] = 46049
$data[6,15] = 5
$data[7,0] = 'Sebastian Olmedo'
$data[7,1] = 'Defensor Central'
$data[7,2] = 54
$data[7,3] = 3650.0264299999999
$data[7,4] = 67.375569999999996
$data[7,5] = 26.440950000000001
$data[7,6] = 167.47
$data[7,7] = 5.19
$data[7,8] = 48.73
$data[7,9] = 0.68
$data[7,10] = 259.81
$data[7,11] = 116.45
$data[7,12] = 6.95
$data[7,13] = 653.16002000000003
$data[7,14] = 46049
$data[7,15] = 5
$data[8,0] = 'Emanuel Beltran'
$data[8,1] = 'Defensor Lateral'
$data[8,2] = 68
$data[8,3] = 4089.9740000000002
$data[8,4] = 60.258040000000001
$data[8,5] = 26.798929999999999
$data[8,6] = 174.94
$data[8,7] = 0
$data[8,8] = 51.26
$data[8,9] = 2.13
$data[8,10] = 394.69
$data[8,11] = 211.51
$data[8,12] = 14.78
$data[8,13] = 833.81
$data[8,14] = 46049
$data[8,15] = 5
$data[9,0] = 'Dario Sarmiento'
$data[9,1] = 'Extremos'
$data[9,2] = 68
$data[9,3] = 3854.5397899999998
$data[9,4] = 56.789360000000002
$data[9,5] = 26.45478
$data[9,6] = 175.63
$data[9,7] = 5.12
$data[9,8] = 45.18
$data[9,9] = 0
$data[9,10] = 260.56000999999998
$data[9,11] = 134.75
$data[9,12] = 12.51
$data[9,13] = 655.91998999999998
$data[9,14] = 46049
$data[9,15] = 5
$data[10,0] = 'Lucas Picech'
$data[10,1] = 'Defensor Central'
$data[10,2] = 54
$data[10,3] = 3847.7670899999998
$data[10,4] = 71.025649999999999
$data[10,5] = 30.234839999999998
$data[10,6] = 183.16
$data[10,7] = 10.93
$data[10,8] = 54.3
$data[10,9] = 1.41
$data[10,10] = 396.42998999999998
$data[10,11] = 139.16999999999999
$data[10,12] = 67.36
$data[10,13] = 677.26999000000001
$data[10,14] = 46049
$data[10,15] = 5
$data[11,0] = 'Alan Bonansea'
$data[11,1] = 'Delantero Central'
$data[11,2] = 54
$data[11,3] = 4580.6323899999998
$data[11,4] = 84.553550000000001
$data[11,5] = 26.740500000000001
$data[11,6] = 186.4
$data[11,7] = 1.28
$data[11,8] = 67.819999999999993
$data[11,9] = 1.84
$data[11,10] = 392.03
$data[11,11] = 280.93000999999998
$data[11,12] = 37.44
$data[11,13] = 872.92997000000003
$data[11,14] = 46049
$data[11,15] = 5
$data[12,0] = 'Baustista Mailler'
$data[12,1] = 'Interno'
$data[12,2] = 54
$data[12,3] = 4716.8524200000002
$data[12,4] = 87.068029999999993
$data[12,5] = 30.834289999999999
$data[12,6] = 197.6
$data[12,7] = 7.29
$data[12,8] = 70.349999999999994
$data[12,9] = 5.33
$data[12,10] = 480.59001000000001
$data[12,11] = 184.69
$data[12,12] = 100.51
$data[12,13] = 945.49998000000005
$data[12,14] = 46049
$data[12,15] = 5
$data[13,0] = 'Ignacio Lagos'
$data[13,1] = 'Extremos'
$data[13,2] = 54
$data[13,3] = 4097.1578099999997
$data[13,4] = 75.629130000000004
$data[13,5] = 28.450710000000001
$data[13,6] = 201.24
$data[13,7] = 8.2200000000000006
$data[13,8] = 63.87
$data[13,9] = 3.19
$data[13,10] = 354.77001999999999
$data[13,11] = 238.01
$data[13,12] = 105.48
$data[13,13] = 827.08996999999999
$data[13,14] = 46049
$data[13,15] = 5
$data[14,0] = 'Ignacio Antonio'
$data[14,1] = 'Volante Central'
$data[14,2] = 54
$data[14,3] = 4143.9860200000003
$data[14,4] = 76.493530000000007
$data[14,5] = 29.864879999999999
$data[14,6] = 201.52
$data[14,7] = 4.99
$data[14,8] = 62.71
$data[14,9] = 2.04
$data[14,10] = 356.92
$data[14,11] = 192.2
$data[14,12] = 48.75
$data[14,13] = 827.08997999999997
$data[14,14] = 46049
$data[14,15] = 5
$data[15,0] = 'Lautaro Gaitan'
$data[15,1] = 'Interno'
$data[15,2] = 67
$data[15,3] = 3612.29376
$data[15,4] = 54.198090000000001
$data[15,5] = 24.3871
$data[15,6] = 203.81
$data[15,7] = 7.48
$data[15,8] = 44.73
$data[15,9] = 1.65
$data[15,10] = 261.27999
$data[15,11] = 86.45
$data[15,12] = 0
$data[15,13] = 538.90000999999995
$data[15,14] = 46049
$data[15,15] = 5
$data[16,0] = 'Facundo Castet'
$data[16,1] = 'Defensor Lateral'
$data[16,2] = 54
$data[16,3] = 4371.0172700000003
$data[16,4] = 80.684280000000001
$data[16,5] = 28.063130000000001
$data[16,6] = 208.38
$data[16,7] = 2.64
$data[16,8] = 57.16
$data[16,9] = 1.4
$data[16,10] = 412.00000999999997
$data[16,11] = 179.53
$data[16,12] = 96.14
$data[16,13] = 852.17002000000002
$data[16,14] = 46049
$data[16,15] = 5
$data[17,0] = 'Matias Muñoz'
$data[17,1] = 'Volante Central'
$data[17,2] = 54
$data[17,3] = 4302.3143600000003
$data[17,4] = 79.4161
$data[17,5] = 27.777979999999999
$data[17,6] = 213.08
$data[17,7] = 4.32
$data[17,8] = 75.22
$data[17,9] = 3.58
$data[17,10] = 412.99000999999998
$data[17,11] = 188.62998999999999
$data[17,12] = 46.83
$data[17,13] = 940.55997000000002
$data[17,14] = 46049
$data[17,15] = 5
$data[18,0] = 'Mauro Peinipil'
$data[18,1] = 'Defensor Lateral'
$data[18,2] = 54
$data[18,3] = 4354.2467699999997
$data[18,4] = 80.374719999999996
$data[18,5] = 28.375229999999998
$data[18,6] = 216.79001
$data[18,7] = 5.55
$data[18,8] = 54.24
$data[18,9] = 1.17
$data[18,10] = 384.26
$data[18,11] = 283.58999999999997
$data[18,12] = 63.33
$data[18,13] = 864.57
$data[18,14] = 46049
$data[18,15] = 5
$data[19,0] = 'Julian Marcioni'
$data[19,1] = 'Extremos'
$data[19,2] = 54
$data[19,3] = 4175.58716
$data[19,4] = 77.076849999999993
$data[19,5] = 29.56073
$data[19,6] = 222.44
$data[19,7] = 7.28
$data[19,8] = 70.77
$data[19,9] = 1.54
$data[19,10] = 436.93999000000002
$data[19,11] = 242.99
$data[19,12] = 75.430000000000007
$data[19,13] = 882.51002000000005
$data[19,14] = 46049
$data[19,15] = 5
$data[20,0] = 'Conrado Ibarra'
$data[20,1] = 'Defensor Lateral'
$data[20,2] = 54
$data[20,3] = 5123.1880499999997
$data[20,4] = 94.568550000000002
$data[20,5] = 31.70234
$data[20,6] = 339.1
$data[20,7] = 27.08
$data[20,8] = 105.52
$data[20,9] = 5.66
$data[20,10] = 591.13996999999995
$data[20,11] = 538.14
$data[20,12] = 275.26
$data[20,13] = 1243.60995
$data[20,14] = 46049
$data[20,15] = 5

$rng = $ws.Range("A470:P490")
$rng.Value = $data

# Apply number formats matching the rest of the sheet:
# Columns D..N and P -> integer-ish numeric format "0"
# Column O -> date format "d-mmm"
$ws.Range("D470:N490").NumberFormat = "0"
$ws.Range("P470:P490").NumberFormat = "0"
$ws.Range("O470:O490").NumberFormat = "d-mmm"
